# EURUSD New Table, Graph Upload
# - Added a new table in the excel
# - Saved the image of the graph

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EURUSD New")
$ws.Activate()

# --- Column widths (approximate best-fit against the engine's quantised
#     character-width model; targets come from the authored OOXML widths) ---
$ws.Columns.Item(1).ColumnWidth = 22.6    # -> ~23.5703125
$ws.Columns.Item(2).ColumnWidth = 10.6    # -> ~11.42578125
$ws.Columns.Item(3).ColumnWidth = 9.76    # -> ~10.7109375
$ws.Columns.Item(4).ColumnWidth = 9.93    # -> ~10.85546875
$ws.Columns.Item(7).ColumnWidth = 9.76    # -> ~10.7109375

# --- New table values for the EURUSD New sheet ---

# Row 2 - ARIMA (D2 keeps the plain/default "Normal" look, no percent style)
$ws.Range("B2").Value = 0.012467
$ws.Range("B2").NumberFormat = "General"
$ws.Range("C2").Value = 0.082795
$ws.Range("C2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("D2").Value = 0.157675
$ws.Range("E2").Value = 0.1505579
$ws.Range("E2").NumberFormat = "General"
$ws.Range("F2").Value = 0.011906
$ws.Range("F2").NumberFormat = "General"
$ws.Range("G2").Value = 0.241715
$ws.Range("G2").NumberFormat = "General"

# Row 3 - LSTM
$ws.Range("B3").Value = -0.096312
$ws.Range("B3").NumberFormat = "General"
$ws.Range("C3").Value = 0.098265
$ws.Range("C3").NumberFormat = "General"
$ws.Range("D3").Value = 0.210969
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = -0.980129
$ws.Range("E3").NumberFormat = "General"
$ws.Range("F3").Value = -0.447452
$ws.Range("F3").NumberFormat = "General"
$ws.Range("G3").Value = -1.445568
$ws.Range("G3").NumberFormat = "General"

# Row 4 - SVM
$ws.Range("B4").Value = -0.096312
$ws.Range("B4").NumberFormat = "General"
$ws.Range("C4").Value = 0.098265
$ws.Range("C4").NumberFormat = "General"
$ws.Range("D4").Value = 0.210969
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = -0.980129
$ws.Range("E4").NumberFormat = "General"
$ws.Range("F4").Value = -0.447452
$ws.Range("F4").NumberFormat = "General"
$ws.Range("G4").Value = -1.445568
$ws.Range("G4").NumberFormat = "General"

# Row 5 - ARIMA-LSTM(1)
$ws.Range("B5").Value = 0.067321
$ws.Range("B5").NumberFormat = "General"
$ws.Range("C5").Value = 0.074702
$ws.Range("C5").NumberFormat = "General"
$ws.Range("D5").Value = 0.064708
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = 0.901185
$ws.Range("E5").NumberFormat = "General"
$ws.Range("F5").Value = 0.952294
$ws.Range("F5").NumberFormat = "General"
$ws.Range("G5").Value = 1.226135
$ws.Range("G5").NumberFormat = "General"

# Row 6 - ARIMA-LSTM(2)
$ws.Range("B6").Value = 0.048557
$ws.Range("B6").NumberFormat = "General"
$ws.Range("C6").Value = 0.082658
$ws.Range("C6").NumberFormat = "General"
$ws.Range("D6").Value = 0.080129
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = 0.587446
$ws.Range("E6").NumberFormat = "General"
$ws.Range("F6").Value = 0.355982
$ws.Range("F6").NumberFormat = "General"
$ws.Range("G6").Value = 0.923387
$ws.Range("G6").NumberFormat = "General"

# Row 7 - ARIMA-SVM(1)
$ws.Range("B7").Value = 0.020189
$ws.Range("B7").NumberFormat = "General"
$ws.Range("C7").Value = 0.082707
$ws.Range("C7").NumberFormat = "General"
$ws.Range("D7").Value = 0.133873
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = 0.244109
$ws.Range("E7").NumberFormat = "General"
$ws.Range("F7").Value = 0.036814
$ws.Range("F7").NumberFormat = "General"
$ws.Range("G7").Value = 0.417182
$ws.Range("G7").NumberFormat = "General"

# Row 8 - ARIMA-SVM(2)
$ws.Range("B8").Value = 0.029375
$ws.Range("B8").NumberFormat = "General"
$ws.Range("C8").Value = 0.099538
$ws.Range("C8").NumberFormat = "General"
$ws.Range("D8").Value = 0.133873
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = 0.29511
$ws.Range("E8").NumberFormat = "General"
$ws.Range("F8").Value = 0.064753
$ws.Range("F8").NumberFormat = "General"
$ws.Range("G8").Value = 0.504343
$ws.Range("G8").NumberFormat = "General"

# Row 9 - GARCH (previously blank)
$ws.Range("B9").Value = -0.037469
$ws.Range("B9").NumberFormat = "General"
$ws.Range("C9").Value = 0.046461
$ws.Range("C9").NumberFormat = "General"
$ws.Range("D9").Value = 0.040779
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = -0.806471
$ws.Range("E9").NumberFormat = "General"
$ws.Range("F9").Value = -0.741016
$ws.Range("F9").NumberFormat = "General"
$ws.Range("G9").Value = -0.974114
$ws.Range("G9").NumberFormat = "General"

# Row 10 - GARCH-LSTM(2) (previously blank)
$ws.Range("B10").Value = 0.051
$ws.Range("B10").NumberFormat = "General"
$ws.Range("C10").Value = 0.05645
$ws.Range("C10").NumberFormat = "General"
$ws.Range("D10").Value = 0.120234
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = 0.903455
$ws.Range("E10").NumberFormat = "General"
$ws.Range("F10").Value = 0.383222
$ws.Range("F10").NumberFormat = "General"
$ws.Range("G10").Value = 1.285767
$ws.Range("G10").NumberFormat = "General"

# --- View state: zoomed in to inspect the new table / graph, selection
#     left where the user clicked after inserting the chart image ---
$excel.ActiveWindow.Zoom = 205
$ws.Range("I10").Select()
